$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "68.064.37"
$ws.Range("E2").Value = "  +0.43%  "

# Row 3
Set-TextValue $ws.Range("D3") "3.793.37"
$ws.Range("E3").Value = "  -0.30%  "

# Row 4
Set-TextValue $ws.Range("D4") "0.998"
$ws.Range("E4").Value = "  -0.27%  "

# Row 5
Set-TextValue $ws.Range("D5") "600.63"
$ws.Range("E5").Value = "  +0.72%  "

# Row 6
Set-TextValue $ws.Range("D6") "165.57"
$ws.Range("E6").Value = "  -0.92%  "

# Row 7
$ws.Range("E7").Value = "  -0.07%  "

# Row 8
$ws.Range("E8").Value = "  -0.71%  "

# Row 9
$ws.Range("E9").Value = "  -1.02%  "

# Row 10
Set-TextValue $ws.Range("D10") "0.452"
$ws.Range("E10").Value = "  +0.42%  "

# Row 11
Set-TextValue $ws.Range("D11") "6.48"
$ws.Range("E11").Value = "  +2.85%  "

# Row 12
Set-TextValue $ws.Range("D12") "0.0000250"
$ws.Range("E12").Value = "  -1.63%  "

# Row 13
Set-TextValue $ws.Range("D13") "35.76"
$ws.Range("E13").Value = "  -0.90%  "

# Row 14
Set-TextValue $ws.Range("D14") "4.427.65"
$ws.Range("E14").Value = "  -0.37%  "

# Row 15
Set-TextValue $ws.Range("D15") "3.792.26"
$ws.Range("E15").Value = "  +0.83%  "

# Row 16
Set-TextValue $ws.Range("D16") "68.036.00"
$ws.Range("E16").Value = "  +0.38%  "

# Row 17
Set-TextValue $ws.Range("D17") "18.43"
$ws.Range("E17").Value = "  -1.06%  "

# Row 18
$ws.Range("E18").Value = "  +1.99%  "

# Row 19
$ws.Range("E19").Value = "  -0.56%  "

# Row 20
Set-TextValue $ws.Range("D20") "461.14"
$ws.Range("E20").Value = "  -0.02%  "

# Row 21
Set-TextValue $ws.Range("D21") "9.71"
$ws.Range("E21").Value = "  -1.73%  "

# Row 22
Set-TextValue $ws.Range("D22") "0.699"
$ws.Range("E22").Value = "  -0.42%  "

# Row 23
$ws.Range("E23").Value = "  -2.88%  "

# Row 24
Set-TextValue $ws.Range("D24") "82.83"
$ws.Range("E24").Value = "  -0.79%  "

# Row 25
Set-TextValue $ws.Range("D25") "12.07"
$ws.Range("E25").Value = "  -0.20%  "

# Row 26
Set-TextValue $ws.Range("D26") "2.11"
$ws.Range("E26").Value = "  +0.31%  "

# Row 27
$ws.Range("E27").Value = "  -0.18%  "

# Row 28
Set-TextValue $ws.Range("D28") "9.99"
$ws.Range("E28").Value = "  -0.18%  "

# Row 29
Set-TextValue $ws.Range("D29") "3.941.79"
$ws.Range("E29").Value = "  -0.24%  "

# Row 30
Set-TextValue $ws.Range("D30") "7.39"
$ws.Range("E30").Value = "  +1.98%  "

# Row 31
Set-TextValue $ws.Range("D31") "2.64"
$ws.Range("E31").Value = "  -5.32%  "

# Row 32
$ws.Range("E32").Value = "  -1.65%  "

# Row 33
Set-TextValue $ws.Range("D33") "29.33"
$ws.Range("E33").Value = "  -1.29%  "

# Row 34
Set-TextValue $ws.Range("D34") "0.999"
$ws.Range("E34").Value = "  -0.07%  "

# Row 35
Set-TextValue $ws.Range("D35") "9.00"
$ws.Range("E35").Value = "  -0.92%  "

# Row 36
Set-TextValue $ws.Range("D36") "0.1000"
$ws.Range("E36").Value = "  -0.04%  "

# Row 37
$ws.Range("E37").Value = "  +0.60%  "

# Row 38
Set-TextValue $ws.Range("D38") "3.27"
$ws.Range("E38").Value = "  -3.10%  "

# Row 39
$ws.Range("B39").Value = "Filecoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D39") "5.79"
$ws.Range("E39").Value = "  -0.13%  "

# Row 40
$ws.Range("B40").Value = "Mantle"
$ws.Range("C40").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws.Range("D40") "0.988"
$ws.Range("E40").Value = "  -0.80%  "

# Row 41
$ws.Range("E41").Value = "  -0.10%  "

# Row 43
Set-TextValue $ws.Range("D43") "0.300"
$ws.Range("E43").Value = "  +0.59%  "

# Row 44
$ws.Range("E44").Value = "  -1.61%  "

# Row 45
Set-TextValue $ws.Range("D45") "43.22"
$ws.Range("E45").Value = "  -1.59%  "

# Row 46
Set-TextValue $ws.Range("D46") "151.59"
$ws.Range("E46").Value = "  +0.69%  "

# Row 47
Set-TextValue $ws.Range("D47") "8.36"
$ws.Range("E47").Value = "  +0.36%  "

# Row 48
Set-TextValue $ws.Range("D48") "1.88"
$ws.Range("E48").Value = "  +2.62%  "

# Row 49
Set-TextValue $ws.Range("D49") "392.57"
$ws.Range("E49").Value = "  +0.18%  "

# Row 50
Set-TextValue $ws.Range("D50") "1.35"
$ws.Range("E50").Value = "  +5.94%  "

# Row 51
$ws.Range("E51").Value = "  +1.13%  "
